$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new localization row: key "test2" -> value "Test 2"
$ws.Range("A5").Value = "test2"
$ws.Range("B5").Value = "Test 2"

# Move the active selection to the newly added row, mirroring Excel's
# behaviour of leaving the cursor on the last-edited cell.
$ws.Range("A5").Select() | Out-Null
